$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "lrs-200 JOB Failure"
$ws.Range("D3").Value = "mcc.prod_INC 3398  ALARM"
$ws.Range("D4").Value = "MAX_RUNTIME"

$ws.Columns.Item(4).ColumnWidth = 19.1

$ws.Range("D3").Select()
